$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (clasificacion) values for specific rows per the diff
$updates = @{
    7 = "Glioma"
    8 = "Glioma"
    11 = "Meningioma"
    13 = "Glioma"
    16 = "Glioma"
    17 = "Glioma"
    25 = "Glioma"
    28 = "Glioma"
    42 = "Meningioma"
    50 = "Glioma"
    56 = "Glioma"
    74 = "Glioma"
    80 = "Glioma"
    82 = "Glioma"
    94 = "Meningioma"
    97 = "Glioma"
    100 = "Glioma"
    104 = "Glioma"
    105 = "Glioma"
    106 = "Glioma"
    113 = "Glioma"
    114 = "Glioma"
    115 = "Meningioma"
    116 = "Pituitary"
    121 = "Glioma"
    122 = "Glioma"
    124 = "Glioma"
    138 = "Glioma"
    152 = "Glioma"
    160 = "Glioma"
    167 = "Meningioma"
    171 = "Pituitary"
    176 = "Glioma"
    191 = "Pituitary"
    201 = "Pituitary"
    206 = "Meningioma"
    210 = "Pituitary"
    212 = "Glioma"
    217 = "Meningioma"
    222 = "Glioma"
    227 = "Glioma"
    239 = "Glioma"
    248 = "Glioma"
    251 = "Meningioma"
    273 = "Glioma"
    278 = "Glioma"
    281 = "Meningioma"
    288 = "Glioma"
    291 = "Meningioma"
    300 = "Glioma"
    310 = "Glioma"
    375 = "Meningioma"
    385 = "Meningioma"
    393 = "Meningioma"
    402 = "Meningioma"
    414 = "Meningioma"
    436 = "Meningioma"
    509 = "Pituitary"
    574 = "Meningioma"
    590 = "Pituitary"
    613 = "Pituitary"
    687 = "Pituitary"
    689 = "Meningioma"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

